# Regenerate the handback status report: refresh the timestamps recorded for
# the most recently processed file (3314774a-5bcf-4947-8455-15fab0d08c42) on
# the Overview sheet and on each per-locale (zh-cn, de-de) sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the .md file's row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-24 16:49:42"

# --- zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-24 16:49:36"
$wsZhCn.Range("K2").Value = "2016-08-24 16:49:54"

# --- de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-24 16:49:42"
$wsDeDe.Range("K2").Value = "2016-08-24 16:50:05"
